# "fix typo in =>" -- Entity Relationships.pptx, slide 12, TextBox 45.
#
# Before:
#   //{ property: p => p=> p.owner }, automatic
#       { property: p .plateNumber },
# After:
#   //{ property: p => p.owner }, automatic
#       { property: p => p.plateNumber },

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(16)
$tr = $shape.TextFrame.TextRange

# --- Fix 1: collapse "p => p=> " into "p => " (drop the duplicated "p=>").
#     Select starting at the run boundary that begins "    //{ property: p "
#     so the replacement merges cleanly into that run's formatting
#     (en-US, Cascadia Mono, green 008000) instead of leaving it split.
$full = $tr.Text
$oldPart1 = "    //{ property: p => p=> "
$idx1 = $full.IndexOf($oldPart1)
$sel1 = $tr.Characters($idx1 + 1, $oldPart1.Length)
$sel1.Text = "    //{ property: p => "

# --- Fix 2: "p .plateNumber" -> "p => p.plateNumber"
# 2a. The lone "." run (own run, de-DE/black) becomes "=> ".
$tr2 = $shape.TextFrame.TextRange
$full2 = $tr2.Text
$oldPart2 = "p .plateNumber"
$idx2 = $full2.IndexOf($oldPart2)
$dotPos = $idx2 + 1 + 2   # 1-based position of the "." character
$selDot = $tr2.Characters($dotPos, 1)
$selDot.Text = "=> "

# 2b. Recolor just the new trailing space to maroon (800000) to match
#     the surrounding punctuation-space runs, splitting it from the "=>".
$tr3 = $shape.TextFrame.TextRange
$selSpace = $tr3.Characters($dotPos + 2, 1)
$selSpace.Font.Color.RGB = 128   # RGB(0x80,0x00,0x00) -> srgbClr 800000

# 2c. "plateNumber" -> "p.plateNumber"
$tr4 = $shape.TextFrame.TextRange
$full4 = $tr4.Text
$idx4 = $full4.IndexOf("plateNumber")
$selPlate = $tr4.Characters($idx4 + 1, "plateNumber".Length)
$selPlate.Text = "p.plateNumber"
